# Database updated with new footprints from EK
# Mark several existing footprint rows as newly "added" (column I, same
# formatting as the existing H "ok"/status column), add a note in I18
# pointing at the source PcbLib the new footprints came from, and for row
# 40 (which already had a note in column I) move the old note into J and
# put the new "added" marker in I.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

$addedRows = @(25, 26, 30, 31, 46, 49, 62, 68, 71, 72, 78, 97, 108)

foreach ($r in $addedRows) {
    $ws.Range("I$r").Value = "added"
    $ws.Range("H$r").Copy()
    $ws.Range("I$r").PasteSpecial($xlPasteFormats)
}

# Row 40 already has a note in I40 ("silk too close to pad"); move that
# note into the new J40 cell (copying I40's formatting) and put the new
# "added" marker in I40.
$ws.Range("J40").Value = "silk too close to pad"
$ws.Range("I40").Copy()
$ws.Range("J40").PasteSpecial($xlPasteFormats)
$ws.Range("I40").Value = "added"

# Note in column I, row 18, referencing the PcbLib the new footprints came
# from (plain/default formatting, no style copied).
$ws.Range("I18").Value = "Footprints\C655_Main_Board_EK.PcbLib"

# New column J holds the (moved) "silk too close to pad" note; give it the
# same width Excel assigned when the column was introduced.
$ws.Columns("J").ColumnWidth = 20.7

$excel.CutCopyMode = $false

# Leave the selection where the editor left it.
[void]$ws.Range("J40").Select()
